# Additional regional data for steel production
# Adds rows 17-26 to the "Fuels" sheet with new fuel/source entries for
# China (CN), Japan (JP) and Russia (RU), plus an IPCC EFDB coke entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write the label ("t=s") cells in the exact order the source
# workbook's shared-string table grew, so new <si> entries line up with the
# authoritative diff (CN coke, CN coking coal, PROXY CN electricity mix,
# JP coking coal, JP steam coal, JP waste plastics, JP coke, RU hard coal,
# IPCC EFDB, IPCC coke, RU natural gas).
$ws.Range("A18").Value = "CN coke"
$ws.Range("A17").Value = "CN coking coal"
$ws.Range("A19").Value = "PROXY CN electricity  mix  (HeEtAl2017)"
$ws.Range("A20").Value = "JP coking coal"
$ws.Range("A21").Value = "JP steam coal"
$ws.Range("A22").Value = "JP waste plastics"
$ws.Range("A23").Value = "JP coke"
$ws.Range("A24").Value = "RU hard coal"
$ws.Range("Q24").Value = "IPCC EFDB"
$ws.Range("A26").Value = "IPCC coke"
$ws.Range("A25").Value = "RU natural gas"

# --- Step 2: fill in the numeric / formula data for each new row.

# Row 17: CN coking coal
$ws.Range("B17").Value = 26.34
$ws.Range("C17").Value = 26.34
$ws.Range("D17").Formula = "=(0.02657*C17)*(44/12)"

# Row 18: CN coke
$ws.Range("B18").Value = 28.435
$ws.Range("C18").Value = 28.435
$ws.Range("D18").Formula = "=(0.02677*C18)*(44/12)"

# Row 19: PROXY CN electricity  mix  (HeEtAl2017)
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Formula = "=(2.9/127.8)*(44/12)"

# Rows 20-23: JP coking coal / JP steam coal / JP waste plastics / JP coke
# (label only, no numeric data in the source edit)

# Row 24: RU hard coal
$ws.Range("B24").Value = 25.16
$ws.Range("C24").Value = 25.16
$ws.Range("D24").Formula = "=93.99/C24"

# Row 25: RU natural gas (label only)

# Row 26: IPCC coke
$ws.Range("B26").Value = 28.2
$ws.Range("C26").Value = 28.2
$ws.Range("D26").Formula = "=C26*29.2*(44/12)/1000"

# --- Step 3: update the view state to match (active cell ends on A26).
$ws.Range("A26").Select()

Write-Output "fuels sheet updated with regional steel-production data"
